$wb = $excel.ActiveWorkbook

# --- Proximity sheet: append two new rows (44, 45) ---
$proximity = $wb.Worksheets.Item("Proximity")

$proxDates = $proximity.Range("A44:A45")
$proxDates.NumberFormat = "@"

$proximity.Cells.Item(44, 1).Value = "2026-02-01"
$proximity.Cells.Item(44, 2).Value = "14:45:11"
$proximity.Cells.Item(44, 3).Value = "14:00"
$proximity.Cells.Item(44, 4).Value = "Living Room Main Door"
$proximity.Cells.Item(44, 5).Value = "ENTER"
$proximity.Cells.Item(44, 6).Value = "User ENTERED Living Room Main Door"

$proximity.Cells.Item(45, 1).Value = "2026-02-01"
$proximity.Cells.Item(45, 2).Value = "14:45:13"
$proximity.Cells.Item(45, 3).Value = "14:00"
$proximity.Cells.Item(45, 4).Value = "Living Room Main Door"
$proximity.Cells.Item(45, 5).Value = "EXIT"
$proximity.Cells.Item(45, 6).Value = "User EXITED Living Room Main Door"

$proxDates.Style = "Normal"

# --- Camera sheet: append one new row (30) ---
$camera = $wb.Worksheets.Item("Camera")

$camDates = $camera.Range("A30")
$camDates.NumberFormat = "@"

$camera.Cells.Item(30, 1).Value = "2026-02-01"
$camera.Cells.Item(30, 2).Value = "14:45:13"
$camera.Cells.Item(30, 3).Value = "14:00"
$camera.Cells.Item(30, 4).Value = "Living Room Main Door"
$camera.Cells.Item(30, 5).Value = "Image Captured"
$camera.Cells.Item(30, 6).Value = "Active"

$camDates.Style = "Normal"
